# Write aircraft to file: Completed. More work on JPADCommander.
# Updates recomputed weight/mass figures on the GLOBAL RESULTS, WING and
# POWER PLANT sheets of the Weights workbook.

$wb = $excel.ActiveWorkbook

# --- GLOBAL RESULTS ---------------------------------------------------
$ws = $wb.Worksheets.Item("GLOBAL RESULTS")

$ws.Range("C6").Value  = 21955.43021694456
$ws.Range("C7").Value  = 21575.43021694456
$ws.Range("C8").Value  = 21296.76731043622
$ws.Range("C12").Value = 3036.2395795835346
$ws.Range("C14").Value = 18919.190637361025
$ws.Range("C15").Value = 18539.190637361025
$ws.Range("C16").Value = 12079.190637361029
$ws.Range("C17").Value = 11849.646991861027
$ws.Range("C18").Value = 11229.22299186103

$ws.Range("C23").Value = 215309.2197369993
$ws.Range("C24").Value = 211582.6927369993
$ws.Range("C25").Value = 208849.94314488926
$ws.Range("C30").Value = 185533.88086387643
$ws.Range("C31").Value = 181807.35386387643
$ws.Range("C32").Value = 118456.39486387651
$ws.Range("C33").Value = 116205.3406727339
$ws.Range("C34").Value = 110121.05965313394

# --- WING ---------------------------------------------------------------
$ws = $wb.Worksheets.Item("WING")

$ws.Range("C9").Value  = 1935.0
$ws.Range("D9").Value  = 11.350884764782073
$ws.Range("C10").Value = 2295.0
$ws.Range("D10").Value = 32.067328441950835
$ws.Range("C13").Value = 1892.2857142857142
$ws.Range("D13").Value = 8.892862280864025

# --- POWER PLANT ---------------------------------------------------------
$ws = $wb.Worksheets.Item("POWER PLANT")

$ws.Range("C2").Value  = 965.2445999999998
$ws.Range("C3").Value  = 1329.1418141999993
$ws.Range("C8").Value  = 482.6223
$ws.Range("C9").Value  = 664.5709070999997
$ws.Range("C12").Value = 482.6223
$ws.Range("C13").Value = 664.5709070999997
